$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuild the LR-pairs table (rows 2-7, cols A-T) with the updated cluster
# naming ("sCs" -> "ECs"/"FAPs" sending clusters) and full set of
# sending/target cluster combinations, following Dr Hou's advice.
$arr = New-Object 'object[,]' 6,20
$arr[0,0] = "ECs"
$arr[0,1] = "Efnb3"
$arr[0,2] = "Ephb1"
$arr[0,3] = "ECs"
$arr[0,4] = 1
$arr[0,5] = 0.3333333333333333
$arr[0,6] = 0.09174599999999999
$arr[0,7] = 0.275238
$arr[0,8] = 0.08724331438250911
$arr[0,9] = 0.08724331438250911
$arr[0,10] = 2
$arr[0,11] = 0.6666666666666666
$arr[0,12] = 1.471191666666667
$arr[0,13] = 4.413575
$arr[0,14] = 0.6447353255635294
$arr[0,15] = 0.6447353255635294
$arr[0,16] = 0.13497595065
$arr[0,17] = 1.21478355585
$arr[0,18] = 0.05624884670164836
$arr[0,19] = 0.05624884670164836
$arr[1,0] = "ECs"
$arr[1,1] = "Efnb3"
$arr[1,2] = "Ephb1"
$arr[1,3] = "sCs"
$arr[1,4] = 1
$arr[1,5] = 0.3333333333333333
$arr[1,6] = 0.09174599999999999
$arr[1,7] = 0.275238
$arr[1,8] = 0.08724331438250911
$arr[1,9] = 0.08724331438250911
$arr[1,10] = 3
$arr[1,11] = 1
$arr[1,12] = 0.8106620000000001
$arr[1,13] = 2.431986
$arr[1,14] = 0.3552646744364706
$arr[1,15] = 0.3552646744364706
$arr[1,16] = 0.07437499585200001
$arr[1,17] = 0.669374962668
$arr[1,18] = 0.03099446768086075
$arr[1,19] = 0.03099446768086075
$arr[2,0] = "FAPs"
$arr[2,1] = "Efnb3"
$arr[2,2] = "Ephb1"
$arr[2,3] = "ECs"
$arr[2,4] = 2
$arr[2,5] = 0.6666666666666666
$arr[2,6] = 0.1498043333333333
$arr[2,7] = 0.449413
$arr[2,8] = 0.1424522763811195
$arr[2,9] = 0.1424522763811195
$arr[2,10] = 2
$arr[2,11] = 0.6666666666666666
$arr[2,12] = 1.471191666666667
$arr[2,13] = 4.413575
$arr[2,14] = 0.6447353255635294
$arr[2,15] = 0.6447353255635294
$arr[2,16] = 0.2203908868305556
$arr[2,17] = 1.983517981475
$arr[2,18] = 0.09184401478984697
$arr[2,19] = 0.09184401478984697
$arr[3,0] = "FAPs"
$arr[3,1] = "Efnb3"
$arr[3,2] = "Ephb1"
$arr[3,3] = "sCs"
$arr[3,4] = 2
$arr[3,5] = 0.6666666666666666
$arr[3,6] = 0.1498043333333333
$arr[3,7] = 0.449413
$arr[3,8] = 0.1424522763811195
$arr[3,9] = 0.1424522763811195
$arr[3,10] = 3
$arr[3,11] = 1
$arr[3,12] = 0.8106620000000001
$arr[3,13] = 2.431986
$arr[3,14] = 0.3552646744364706
$arr[3,15] = 0.3552646744364706
$arr[3,16] = 0.1214406804686667
$arr[3,17] = 1.092966124218
$arr[3,18] = 0.05060826159127256
$arr[3,19] = 0.05060826159127256
$arr[4,0] = "sCs"
$arr[4,1] = "Efnb3"
$arr[4,2] = "Ephb1"
$arr[4,3] = "ECs"
$arr[4,4] = 3
$arr[4,5] = 1
$arr[4,6] = 0.8100603333333333
$arr[4,7] = 2.430181
$arr[4,8] = 0.7703044092363713
$arr[4,9] = 0.7703044092363713
$arr[4,10] = 2
$arr[4,11] = 0.6666666666666666
$arr[4,12] = 1.471191666666667
$arr[4,13] = 4.413575
$arr[4,14] = 0.6447353255635294
$arr[4,15] = 0.6447353255635294
$arr[4,16] = 1.191754011897222
$arr[4,17] = 10.725786107075
$arr[4,18] = 0.4966424640720341
$arr[4,19] = 0.4966424640720341
$arr[5,0] = "sCs"
$arr[5,1] = "Efnb3"
$arr[5,2] = "Ephb1"
$arr[5,3] = "sCs"
$arr[5,4] = 3
$arr[5,5] = 1
$arr[5,6] = 0.8100603333333333
$arr[5,7] = 2.430181
$arr[5,8] = 0.7703044092363713
$arr[5,9] = 0.7703044092363713
$arr[5,10] = 3
$arr[5,11] = 1
$arr[5,12] = 0.8106620000000001
$arr[5,13] = 2.431986
$arr[5,14] = 0.3552646744364706
$arr[5,15] = 0.3552646744364706
$arr[5,16] = 0.6566851299406667
$arr[5,17] = 5.910166169466
$arr[5,18] = 0.2736619451643373
$arr[5,19] = 0.2736619451643373
$ws.Range("A2:T7").Value = $arr